# Updates the Price (D) and Volume(1h) (E) columns for the cryptos list
# to the latest scraped values, per the GitHub Actions refresh commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "25.728.95" }
    @{ Cell = "E2"; Value = "  -4.21%  " }
    @{ Cell = "D3"; Value = "1.816.95" }
    @{ Cell = "E3"; Value = "  -3.04%  " }
    @{ Cell = "D4"; Value = "1.001" }
    @{ Cell = "D5"; Value = "277.29" }
    @{ Cell = "E5"; Value = "  -8.01%  " }
    @{ Cell = "E6"; Value = "  -0.13%  " }
    @{ Cell = "D7"; Value = "0.5089" }
    @{ Cell = "E7"; Value = "  -4.75%  " }
    @{ Cell = "D8"; Value = "0.3525" }
    @{ Cell = "E8"; Value = "  -6.03%  " }
    @{ Cell = "D9"; Value = "44.37" }
    @{ Cell = "E9"; Value = "  -2.46%  " }
    @{ Cell = "D10"; Value = "0.06668" }
    @{ Cell = "E10"; Value = "  -7.18%  " }
    @{ Cell = "D11"; Value = "20.05" }
    @{ Cell = "E11"; Value = "  -7.04%  " }
    @{ Cell = "D12"; Value = "0.8270" }
    @{ Cell = "E12"; Value = "  -7.03%  " }
    @{ Cell = "D13"; Value = "0.07888" }
    @{ Cell = "E13"; Value = "  -3.19%  " }
    @{ Cell = "D14"; Value = "1.818.40" }
    @{ Cell = "E14"; Value = "  -2.98%  " }
    @{ Cell = "D15"; Value = "5.069" }
    @{ Cell = "D16"; Value = "87.55" }
    @{ Cell = "E16"; Value = "  -6.11%  " }
    @{ Cell = "D17"; Value = "1.000" }
    @{ Cell = "E17"; Value = "  -0.13%  " }
    @{ Cell = "D18"; Value = "14.09" }
    @{ Cell = "E18"; Value = "  -5.06%  " }
    @{ Cell = "D19"; Value = "0.000008027" }
    @{ Cell = "E19"; Value = "  -5.95%  " }
    @{ Cell = "E20"; Value = "  -0.09%  " }
    @{ Cell = "D21"; Value = "25.777.67" }
    @{ Cell = "E21"; Value = "  -4.16%  " }
    @{ Cell = "D22"; Value = "4.743" }
    @{ Cell = "E22"; Value = "  -4.81%  " }
    @{ Cell = "D23"; Value = "10.00" }
    @{ Cell = "E24"; Value = "  -4.78%  " }
    @{ Cell = "D25"; Value = "142.11" }
    @{ Cell = "E25"; Value = "  -2.76%  " }
    @{ Cell = "D26"; Value = "2.203" }
    @{ Cell = "E26"; Value = "  -4.09%  " }
    @{ Cell = "D27"; Value = "1.673" }
    @{ Cell = "E27"; Value = "  -3.36%  " }
    @{ Cell = "D28"; Value = "17.10" }
    @{ Cell = "E28"; Value = "  -5.50%  " }
    @{ Cell = "D29"; Value = "109.47" }
    @{ Cell = "E29"; Value = "  -3.98%  " }
    @{ Cell = "E30"; Value = "  -8.31%  " }
    @{ Cell = "D31"; Value = "4.232" }
    @{ Cell = "E31"; Value = "  -8.19%  " }
    @{ Cell = "D32"; Value = "0.08772" }
    @{ Cell = "E32"; Value = "  -4.06%  " }
    @{ Cell = "D33"; Value = "0.04882" }
    @{ Cell = "E33"; Value = "  -2.61%  " }
    @{ Cell = "D34"; Value = "0.7268" }
    @{ Cell = "D35"; Value = "1.136" }
    @{ Cell = "E35"; Value = "  -3.20%  " }
    @{ Cell = "E36"; Value = "  -2.75%  " }
    @{ Cell = "D37"; Value = "3.132" }
    @{ Cell = "E37"; Value = "  -2.55%  " }
    @{ Cell = "D38"; Value = "2.372" }
    @{ Cell = "E38"; Value = "  -9.18%  " }
    @{ Cell = "D39"; Value = "0.01851" }
    @{ Cell = "E39"; Value = "  -5.28%  " }
    @{ Cell = "D40"; Value = "0.5174" }
    @{ Cell = "E40"; Value = "  -13.99%  " }
    @{ Cell = "D41"; Value = "0.9648" }
    @{ Cell = "E41"; Value = "  -9.78%  " }
    @{ Cell = "E42"; Value = "  -6.17%  " }
    @{ Cell = "D43"; Value = "110.38" }
    @{ Cell = "E43"; Value = "  -4.06%  " }
    @{ Cell = "D44"; Value = "8.023" }
    @{ Cell = "E44"; Value = "  -10.03%  " }
    @{ Cell = "E45"; Value = "  -0.09%  " }
    @{ Cell = "D46"; Value = "0.4550" }
    @{ Cell = "E46"; Value = "  -10.56%  " }
    @{ Cell = "D47"; Value = "0.1364" }
    @{ Cell = "E47"; Value = "  -8.56%  " }
    @{ Cell = "D48"; Value = "36.49" }
    @{ Cell = "E48"; Value = "  -3.19%  " }
    @{ Cell = "D49"; Value = "9.264" }
    @{ Cell = "E49"; Value = "  -6.82%  " }
    @{ Cell = "D50"; Value = "1.501" }
    @{ Cell = "E50"; Value = "  -8.14%  " }
    @{ Cell = "D51"; Value = "0.05843" }
    @{ Cell = "E51"; Value = "  -3.44%  " }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    # The Price column contains values that look numeric (e.g. "1.000", "10.00").
    # Force them to be stored as text so formatting/precision (trailing zeros,
    # thousands-style dots) is preserved exactly, matching the scraped strings.
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}
